$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 19.2
$ws.Range("C3").Value = 20.5
$ws.Range("B4").Value = 21.5
$ws.Range("C4").Value = 20.5
$ws.Range("C5").Value = 19.5
$ws.Range("C17").Value = 13
$ws.Range("C22").Value = 16.8
